# Automatically re-generate list and index
#
# Every "review date" string stored as text in the form YYYY-MM-20
# (day-of-month 20) is bumped one day forward to YYYY-MM-21. Dates that
# don't end in "-20" (e.g. 2024-06-02, 2024-06-26, 2024-07-25, ...) are
# left untouched.
#
# The cells hold plain text (not real Excel dates), so we briefly force a
# Text number format while writing the new literal string - otherwise
# Excel's automatic type detection would silently turn "2011-03-21" into
# a date serial number - and then restore the cell's original ("Normal")
# style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1

$pattern = '^(\d{4}-\d{2}-)20$'

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = [string]$cell.Value2

        if ($text -match $pattern) {
            $newText = $text.Substring(0, $text.Length - 2) + "21"

            $cell.NumberFormat = "@"
            $cell.Value = $newText
            $cell.Style = "Normal"
        }
    }
}
